$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 156, shifting existing rows 156-229 down to 157-230.
$ws.Rows("156:156").Insert()

# Populate the newly inserted row 156 with the new market-price record.
$ws.Range("A156").Value = 10
$ws.Range("B156").Value = "Vega Modelo de Temuco"
$ws.Range("C156").Value = "La Araucanía"
$ws.Range("D156").Value = 44523
$ws.Range("E156").Value = 9
$ws.Range("F156").Value = 100114013
$ws.Range("G156").Value = "Zanahoria"
$ws.Range("H156").Value = "Sin especificar"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 170
$ws.Range("K156").Value = 7000
$ws.Range("L156").Value = 8000
$ws.Range("M156").Value = 7471
$ws.Range("N156").Value = '$/saco 20 kilos'
$ws.Range("O156").Value = "Provincia del Elquí"
$ws.Range("P156").Value = 374
$ws.Range("Q156").Value = 20
$ws.Range("R156").Value = "Hortaliza"
